# Append, after the existing content, an empty paragraph followed by a
# paragraph of underlined text "Aloha asdmasdasdasd".
#
# We build the new paragraphs via a raw WordprocessingML fragment and
# insert it at the very end of the document's main story so that the
# resulting OOXML shape matches exactly what Word itself would produce
# (an empty self-closed <w:p/> then a <w:p> whose paragraph mark carries
# the underline formatting while the run text itself stays unformatted).

$d = $word.ActiveDocument

$end = $d.Content
$end.Collapse(0)   # wdCollapseEnd -> collapse to a zero-length range at the end of the document

$fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' +
                    '<w:p/>' +
                    '<w:p>' +
                        '<w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
                        '<w:r><w:t>Aloha asdmasdasdasd</w:t></w:r>' +
                    '</w:p>' +
                '</w:body>' +
            '</w:document>' +
        '</pkg:xmlData>' +
    '</pkg:part>' +
'</pkg:package>'

$end.InsertXML($fragment)
